$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 16:54:15"
$wsZh.Range("H2").Value = "2016-03-19 16:54:32"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 16:54:18"
$wsDe.Range("H2").Value = "2016-03-19 16:54:37"
